# Complete summary of statistics
# Adds a new row (row 15) to the ANOVA summary table on Sheet1, continuing
# the "No-Cursor VS No-Cursor" comparison family started in row 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for row 15 (Time / Experiment / Time*Experiment / Comparison) ---
# Comparison label entered first, then the three statistics, to match the
# shared-string insertion order of the original edit.
$ws.Range("D15").Value = "No-Cursor VS No-Cursor -No- Cursors"
$ws.Range("A15").Value = "F(3,120) = 106.74, p < .001"
$ws.Range("B15").Value = "F(1,40) = 3.01, p = 0.09"
$ws.Range("C15").Value = "F(3,120) = 1.43, p = 0.239"

# --- Match the formatting of the row above (same comparison family) ---
# A15/C15 share the "light" shade, B15 the "dark" shade (like A14/B14/C14/D14).
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

# --- Column D needs to widen to fit the new, longer comparison label ---
# (matches the "best fit" width Excel computes for the new, longer text)
$ws.Columns.Item(4).ColumnWidth = 33.6

# --- Restore the selection the author left the sheet in ---
[void]$ws.Range("D28").Select()
